$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.773.81'
$ws.Range('E2').Value = '  +0.64%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.918.00'
$ws.Range('E3').Value = '  +1.61%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9923'
$ws.Range('E4').Value = '  -0.95%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '250.73'
$ws.Range('E5').Value = '  +2.42%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5982'
$ws.Range('E6').Value = '  +26.82%  '

$ws.Range('E7').Value = '  -0.90%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3042'
$ws.Range('E8').Value = '  +4.42%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '24.36'
$ws.Range('E9').Value = '  +8.54%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06646'
$ws.Range('E10').Value = '  +2.31%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7939'
$ws.Range('E11').Value = '  +7.56%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '101.84'
$ws.Range('E12').Value = '  +6.02%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07888'
$ws.Range('E13').Value = '  +1.59%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.898.28'
$ws.Range('E14').Value = '  +0.57%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.323'
$ws.Range('E15').Value = '  +2.53%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '285.60'
$ws.Range('E16').Value = '  +0.68%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.676.67'
$ws.Range('E17').Value = '  +0.07%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.51'
$ws.Range('E18').Value = '  +3.46%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007611'
$ws.Range('E19').Value = '  +1.45%  '

$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.461'
$ws.Range('E20').Value = '  +3.64%  '

$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9928'
$ws.Range('E21').Value = '  -0.83%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.144.68'
$ws.Range('E22').Value = '  +0.86%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9914'
$ws.Range('E23').Value = '  -1.01%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.564'
$ws.Range('E24').Value = '  +4.84%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.279'
$ws.Range('E25').Value = '  +1.22%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.95'
$ws.Range('E26').Value = '  -0.24%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.45'
$ws.Range('E27').Value = '  +3.13%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.958'
$ws.Range('E28').Value = '  +2.84%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1069'
$ws.Range('E29').Value = '  +9.65%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.341'
$ws.Range('E30').Value = '  -0.55%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.522'
$ws.Range('E31').Value = '  +3.06%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.414'
$ws.Range('E32').Value = '  +2.69%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.275'
$ws.Range('E33').Value = '  +3.41%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04930'
$ws.Range('E34').Value = '  +1.05%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.155'
$ws.Range('E35').Value = '  +2.30%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7137'
$ws.Range('E36').Value = '  +2.94%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.773'
$ws.Range('E37').Value = '  +2.37%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01930'

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.904'
$ws.Range('E39').Value = '  +2.37%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '76.96'
$ws.Range('E40').Value = '  +1.94%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.378'
$ws.Range('E41').Value = '  +2.84%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4389'
$ws.Range('E42').Value = '  +2.91%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.004'
$ws.Range('E43').Value = '  -0.23%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8452'
$ws.Range('E44').Value = '  +2.44%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9933'
$ws.Range('E45').Value = '  -0.80%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.05'
$ws.Range('E46').Value = '  +5.54%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.30'
$ws.Range('E47').Value = '  -0.15%  '

$ws.Range('E48').Value = '  +3.24%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '35.54'
$ws.Range('E49').Value = '  +0.33%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4058'
$ws.Range('E50').Value = '  +2.83%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '905.53'
$ws.Range('E51').Value = '  -0.62%  '
